$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows where "dplyr" used to be (row 46), pushing it down to row 48.
$ws.Rows("46:47").Insert()

# New packages added above the (now relocated) "dplyr" row.
$ws.Range("A46").Value = "visdat"
$ws.Range("A47").Value = "tidymodels"

# The relocated "dplyr" cell picks up a distinct (explicitly-applied) style.
$ws.Range("A48").Style = "Normal"

# Reflect the author's final cursor/scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("D47").Select() | Out-Null
